# Reg TCs for 811, 2805, 2648 stories
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the new "ExpectedFilenames" values for the second scenario block
# (rows 7-10) and add a brand new row 11 for an additional filename.
$ws.Range("J7").Value  = "StandardExcelReport-Test_Automation_1-Clinical-2023_"
$ws.Range("J8").Value  = "ExcelReport-Test_Automation_1-Clinical-"
$ws.Range("J9").Value  = "WordReport-Test_Automation_1-Clinical-"
$ws.Range("J10").Value = "CompleteExcelReport-Test_NonOncology_Automation_3-Clinical-2023_"
$ws.Range("J11").Value = "StandardExcelReport-Test_NonOncology_Automation_3-Clinical-2023_"

# Move the active selection the way the author left it when saving.
$ws.Range("J15").Select()
